$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "CPU vs GPU" summary table (rows 2-15): recomputed summation values ---
# Column C = CPU time (ms), Column D = GPU time (ms)

$ws.Cells.Item(2, 4).Value = 672
$ws.Cells.Item(3, 4).Value = 83
$ws.Cells.Item(4, 4).Value = 89
$ws.Cells.Item(5, 4).Value = 74
$ws.Cells.Item(7, 4).Value = 97
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(8, 4).Value = 76
$ws.Cells.Item(9, 3).Value = 1
$ws.Cells.Item(9, 4).Value = 97
$ws.Cells.Item(10, 4).Value = 110
$ws.Cells.Item(11, 4).Value = 105
$ws.Cells.Item(12, 3).Value = 4
$ws.Cells.Item(12, 4).Value = 135
$ws.Cells.Item(13, 3).Value = 10
$ws.Cells.Item(13, 4).Value = 196
$ws.Cells.Item(14, 4).Value = 323
$ws.Cells.Item(15, 3).Value = 208
$ws.Cells.Item(15, 4).Value = 936

# --- Shrink the shared-formula group for the third "2^i" table (rows 92:115) ---
$ws.Range("B92:B115").Formula = "=2^A92"

# --- Optimize the summation: turn the fourth table's (rows 122:143) individually
#     entered formulas into one shared formula, like the other tables on the sheet ---
$ws.Range("B122:B143").Formula = "=2^A122"

# --- Restore the sheet view: reset scroll position, update zoom, move selection ---
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$win.Zoom = 85
$ws.Range("D16").Select()
